# Update CDA Logical model for ST.r2b
# - bump Version / Date metadata values
# - insert a new "Jurisdiction" property row (blank value) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version value (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for the "Jurisdiction" property
$ws.Rows.Item(11).Insert()

# Match formatting of the surrounding data rows
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

Write-Host "Edit complete"
